# Append new DRS rows for matches 42-47 (data updated till 47 matches)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(42,"KKR","PBKS",1,"KKR","PBKS",9,"KKR","YC Barde","YCB","Wicket","Out","Not Out","SP Narine","RD Chahar","Successful","No"),
    @(42,"KKR","PBKS",1,"KKR","PBKS",20,"KKR","AK Chaudhary","AKC","Wide","Not Called","Not Called","Ramandeep Singh","HV Patel","Unsuccessful","No"),
    @(42,"KKR","PBKS",1,"KKR","PBKS",20,"PBKS","AK Chaudhary","AKC","Wide","Called","Not Called","VR Iyer","HV Patel","Successful","No"),
    @(42,"KKR","PBKS",2,"PBKS","KKR",17,"KKR","YC Barde","YCB","NoBall","Called","Not Called","Shashank Singh","PVD Chameera","Successful","No"),
    @(42,"KKR","PBKS",2,"PBKS","KKR",18,"PBKS","AK Chaudhary","AKC","Wide","Not Called","Called","Shashank Singh","Harshit Rana","Successful","No"),
    @(43,"DC","MI",1,"DC","MI",15,"MI","Navdeep Singh","NS","Wicket","Not Out","Not Out","RR Pant","PP Chawla","Unsuccessful","No"),
    @(43,"DC","MI",2,"MI","DC",11,"DC","NA Patwardhan","NAP","Wicket","Not Out","Not Out","HH Pandya","Rasikh Salam","Unsuccessful","No"),
    @(43,"DC","MI",2,"MI","DC",18,"MI","Navdeep Singh","NS","Wicket","Out","Out","TH David","Mukesh Kumar","Unsuccessful","No"),
    @(43,"DC","MI",2,"MI","DC",19,"MI","NA Patwardhan","NAP","Wide","Not Called","Not Called","Mohammad Nabi","Rasikh Salam","Unsuccessful","No"),
    @(44,"LSG","RR",1,"LSG","RR",19,"RR","MA Gough","MAG","Wide","Called","Called","A Badoni","TA Boult","Unsuccessful","No"),
    @(45,"GT","RCB",1,"GT","RCB",20,"GT","Nitin Menon","NM","NoBall","Not Called","Not Called","DA Miller","Yash Dayal","Unsuccessful","No"),
    @(45,"GT","RCB",2,"RCB","GT",7,"GT","VK Sharma","VKS","Wicket","Not Out","Not Out","WG Jacks","Noor Ahmad","Unsuccessful","No"),
    @(45,"GT","RCB",2,"RCB","GT",10,"RCB","Nitin Menon","NM","Wicket","Out","Not Out","WG Jacks","Rashid Khan","Successful","No"),
    @(46,"CSK","SRH",1,"CSK","SRH",2,"SRH","R Pandit","RP","Wide","Called","Not Called","AM Rahane","K Nitish Kumar Reddy","Successful","No"),
    @(46,"CSK","SRH",1,"CSK","SRH",17,"CSK","MV Saidharshan Kumar","MVSK","Wide","Not Called","Called","S Dube","T Natarajan","Successful","No"),
    @(46,"CSK","SRH",2,"SRH","CSK",17,"SRH","MV Saidharshan Kumar","MVSK","Wide","Not Called","Not Called","PJ Cummins","SN Thakur","Unsuccessful","No"),
    @(46,"CSK","SRH",2,"SRH","CSK",18,"SRH","R Pandit","RP","Wide","Not Called","Not Called","PJ Cummins","TU Deshpande","Unsuccessful","No"),
    @(47,"KKR","DC",1,"DC","KKR",2,"KKR","Tapan Sharma","TS","Wicket","Not Out","Out","PP Shaw","VG Arora","Successful","No"),
    @(47,"KKR","DC",1,"DC","KKR",15,"KKR","Navdeep Singh","NS","Wicket","Not Out","Not Out","Rasikh Salam","CV Varun","Unsuccessful","No"),
    @(47,"KKR","DC",1,"DC","KKR",18,"KKR","Tapan Sharma","TS","Wicket","Not Out","Not Out","Kuldeep Yadav","VG Arora","Unsuccessful","No"),
    @(47,"KKR","DC",1,"DC","KKR",20,"DC","Tapan Sharma","TS","Wide","Not Called","Not Called","Kuldeep Yadav","AD Russell","Unsuccessful","No"),
    @(47,"KKR","DC",2,"KKR","DC",2,"DC","Tapan Sharma","TS","Wide","Called","Not Called","PD Salt","KK Ahmed","Successful","No")
)

$nRows = $rows.Count
$nCols = $rows[0].Count

$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

# New data starts right after the last existing row (170)
$startRow = 171
$endRow = $startRow + $nRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,$nCols))
$rng.Value = $arr

# Column O (Bowler) auto-widened slightly by Excel after the longer new names were added
$ws.Columns.Item(15).ColumnWidth = 17.5

# Leave the selection where the author ended up after the edit
$ws.Range("P188").Select() | Out-Null

Write-Host "Added $nRows rows ($startRow-$endRow)."
